$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 536, pushing existing rows 536:574 down to 537:575.
$ws.Rows.Item(536).Insert()

# Populate the newly-inserted row 536 with the new data record.
$ws.Range("A536").Value = 10
$ws.Range("B536").Value = "Vega Modelo de Temuco"
$ws.Range("C536").Value = "La Araucanía"
$ws.Range("D536").Value = 44931
$ws.Range("E536").Value = 9
$ws.Range("F536").Value = 100112023
$ws.Range("G536").Value = "Brócoli"
$ws.Range("H536").Value = "Sin especificar"
$ws.Range("I536").Value = "Primera"
$ws.Range("J536").Value = 800
$ws.Range("K536").Value = 1200
$ws.Range("L536").Value = 1200
$ws.Range("M536").Value = 1200
$ws.Range("N536").Value = "$/unidad"
$ws.Range("O536").Value = "Provincia de Cautín"
$ws.Range("P536").Value = 1200
$ws.Range("Q536").Value = 1
$ws.Range("R536").Value = "Hortaliza"

# Match the date-format style used by the rest of column D.
$ws.Range("D536").NumberFormat = $ws.Range("D537").NumberFormat
